$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.02;
    "C2" = 1.028178317714243;
    "D2" = 1.032125997470701;
    "E2" = 1.028177753935987;
    "F2" = 1.038536709080203;
    "I2" = 1.034782100449942;
    "J2" = 1.033332078246427;
    "K2" = 1.034932128703562;
    "L2" = 1.030995328357599;
    "M2" = 1.041324461438989;
    "N2" = 1.005712725503983;
    "B3" = 1.02;
    "C3" = 1.02901428895391;
    "D3" = 1.032734477795152;
    "E3" = 1.028884844064936;
    "F3" = 1.039630400522287;
    "I3" = 1.034973945342576;
    "J3" = 1.033809130661187;
    "K3" = 1.035349852419388;
    "L3" = 1.031510571376367;
    "M3" = 1.04222743734959;
    "B4" = 1.02;
    "C4" = 1.029555748876456;
    "D4" = 1.03312863166901;
    "E4" = 1.029343219136517;
    "F4" = 1.040339000620872;
    "I4" = 1.035097165090164;
    "J4" = 1.034117694666493;
    "K4" = 1.035619879516503;
    "L4" = 1.031844135539274;
    "M4" = 1.042812041727811;
    "B5" = 1.02;
    "C5" = 1.029783504251742;
    "D5" = 1.033294434478396;
    "E5" = 1.029536119598433;
    "F5" = 1.040637112796076;
    "I5" = 1.03514874661287;
    "J5" = 1.034247384981626;
    "K5" = 1.03573333367174;
    "L5" = 1.031984404809889;
    "M5" = 1.043057884838283;
    "B6" = 1.02;
    "C6" = 1.029821752707136;
    "D6" = 1.033322279348732;
    "E6" = 1.029568520090947;
    "F6" = 1.04068717982518;
    "I6" = 1.035157394455924;
    "J6" = 1.03426915876628;
    "K6" = 1.035752379267388;
    "L6" = 1.032007958883724;
    "M6" = 1.043099167377257;
    "B7" = 1.02;
    "C7" = 1.029558791663072;
    "D7" = 1.033130846741254;
    "E7" = 1.02934579589976;
    "F7" = 1.040342983159912;
    "I7" = 1.035097855190414;
    "J7" = 1.034119427714024;
    "K7" = 1.035621395754138;
    "L7" = 1.031846009672801;
    "M7" = 1.042815326400139;
    "B8" = 1.02;
    "C8" = 1.028460727389323;
    "D8" = 1.03233154701322;
    "E8" = 1.028416543349912;
    "F8" = 1.038906138591873;
    "I8" = 1.0348471247867;
    "J8" = 1.033493324713046;
    "K8" = 1.035073355260766;
    "L8" = 1.031169421751041;
    "M8" = 1.041629559904111;
    "B9" = 1.02;
    "C9" = 1.026529924223539;
    "D9" = 1.030926413668264;
    "E9" = 1.02678559067457;
    "F9" = 1.036381230288489;
    "I9" = 1.034398311009536;
    "J9" = 1.032389173337054;
    "K9" = 1.034105631452347;
    "L9" = 1.029978522967105;
    "M9" = 1.039542564945297;
    "B10" = 1.02;
    "C10" = 1.025245580053754;
    "D10" = 1.029991993345788;
    "E10" = 1.025702762192786;
    "F10" = 1.034702709201028;
    "I10" = 1.034094434196457;
    "J10" = 1.031652545749159;
    "K10" = 1.033459196058225;
    "L10" = 1.029185559573804;
    "M10" = 1.038152949254817;
    "B11" = 1.02;
    "C11" = 1.024690140174737;
    "D11" = 1.029587952525911;
    "E11" = 1.025234966951198;
    "F11" = 1.033977027508453;
    "I11" = 1.033961752456589;
    "J11" = 1.031333466273653;
    "K11" = 1.033178990516856;
    "L11" = 1.028842442486907;
    "M11" = 1.037551646796448;
    "B12" = 1.02;
    "C12" = 1.024483930066112;
    "D12" = 1.029437960895254;
    "E12" = 1.025061370374549;
    "F12" = 1.033707647347245;
    "I12" = 1.033912303700407;
    "J12" = 1.031214929750629;
    "K12" = 1.033074866469825;
    "L12" = 1.028715030973268;
    "M12" = 1.03732835867328;
    "B13" = 1.02;
    "C13" = 1.024528158060959;
    "D13" = 1.029470130641466;
    "E13" = 1.025098600008619;
    "F13" = 1.033765422586991;
    "I13" = 1.033922918077061;
    "J13" = 1.031240356953224;
    "K13" = 1.03309720337089;
    "L13" = 1.028742359461548;
    "M13" = 1.037376251898779;
    "B14" = 1.02;
    "C14" = 1.024673092632972;
    "D14" = 1.029575552384575;
    "E14" = 1.025220614061066;
    "F14" = 1.033954756981479;
    "I14" = 1.033957668366549;
    "J14" = 1.031323668332007;
    "K14" = 1.033170384466511;
    "L14" = 1.028831909840555;
    "M14" = 1.037533188453641;
    "B15" = 1.02;
    "C15" = 1.024762405560379;
    "D15" = 1.029640517791347;
    "E15" = 1.025295812647836;
    "F15" = 1.034071434742511;
    "I15" = 1.033979057339732;
    "J15" = 1.031374997107281;
    "K15" = 1.033215468057444;
    "L15" = 1.028887089782868;
    "M15" = 1.037629890533211;
    "B16" = 1.02;
    "C16" = 1.02528245741696;
    "D16" = 1.030018820334892;
    "E16" = 1.025733831046899;
    "F16" = 1.034750894152351;
    "I16" = 1.034103216697102;
    "J16" = 1.031673719673311;
    "K16" = 1.03347778625676;
    "L16" = 1.029208336310907;
    "M16" = 1.038192864426242;
    "B17" = 1.02;
    "C17" = 1.025608857815465;
    "D17" = 1.030256273073561;
    "E17" = 1.02600887782955;
    "F17" = 1.035177403856394;
    "I17" = 1.034180804197168;
    "J17" = 1.031861070440594;
    "K17" = 1.033642253361207;
    "L17" = 1.029409911218112;
    "M17" = 1.038546113446029;
    "B18" = 1.02;
    "C18" = 1.025799308091089;
    "D18" = 1.030394830066339;
    "E18" = 1.02616941173359;
    "F18" = 1.035426288636412;
    "I18" = 1.034225953396978;
    "J18" = 1.031970337819205;
    "K18" = 1.033738155731168;
    "L18" = 1.029527509652299;
    "M18" = 1.038752197161395;
    "B19" = 1.02;
    "C19" = 1.025864257941591;
    "D19" = 1.030442083671801;
    "E19" = 1.026224167197049;
    "F19" = 1.035511170362081;
    "I19" = 1.034241330063624;
    "J19" = 1.032007593247427;
    "K19" = 1.033770851093332;
    "L19" = 1.029567611578213;
    "M19" = 1.038822473066704;
    "B20" = 1.02;
    "C20" = 1.025573831258954;
    "D20" = 1.030230790959424;
    "E20" = 1.02597935716753;
    "F20" = 1.035131632133879;
    "I20" = 1.034172490781207;
    "J20" = 1.031840970622418;
    "K20" = 1.033624610538689;
    "L20" = 1.029388281729795;
    "M20" = 1.038508209074278;
    "B21" = 1.02;
    "C21" = 1.024630410104439;
    "D21" = 1.029544505891175;
    "E21" = 1.025184679439715;
    "F21" = 1.033898998013657;
    "I21" = 1.033947439818483;
    "J21" = 1.031299135646097;
    "K21" = 1.033148835642847;
    "L21" = 1.028805538442511;
    "M21" = 1.037486972821736;
    "B22" = 1.02;
    "C22" = 1.024037851927992;
    "D22" = 1.029113516059351;
    "E22" = 1.024685980236798;
    "F22" = 1.03312497785472;
    "I22" = 1.03380498785816;
    "J22" = 1.030958369294332;
    "K22" = 1.032849447470368;
    "L22" = 1.028439361873809;
    "M22" = 1.036845242355699;
    "B23" = 1.02;
    "C23" = 1.02435192012003;
    "D23" = 1.029341943544313;
    "E23" = 1.024950259877632;
    "F23" = 1.03353520690147;
    "I23" = 1.033880594525159;
    "J23" = 1.031139024465182;
    "K23" = 1.033008182114949;
    "L23" = 1.028633458007143;
    "M23" = 1.037185401322252;
    "B24" = 1.02;
    "C24" = 1.025589658046985;
    "D24" = 1.030242305060725;
    "E24" = 1.02599269596357;
    "F24" = 1.035152314068625;
    "I24" = 1.03417624758447;
    "J24" = 1.031850052899649;
    "K24" = 1.033632582658969;
    "L24" = 1.029398055092988;
    "M24" = 1.038525336306147;
    "B25" = 1.02;
    "C25" = 1.027028585523107;
    "D25" = 1.031289269519513;
    "E25" = 1.027206449678704;
    "F25" = 1.037033145716089;
    "I25" = 1.034515165419162;
    "J25" = 1.032674720419604;
    "K25" = 1.034356042552321;
    "L25" = 1.030286233240774;
    "M25" = 1.04008180460958
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"
